# ---------------------------------------------------------------------------
# SPM changes 2012 SM.docx - apply commit edits via Word COM-interop
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: remove yellow highlighting from the
# "age_size_distribution normal_by_length " paragraph (run + paragraph mark).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("age_size_distribution normal_by_length")
if ($found1) {
    $p1 = $r1.Paragraphs(1)
    $pr1 = $p1.Range
    $xml1 = "<w:p $wNs>" +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r><w:t>age_size_</w:t></w:r>' +
              '<w:r><w:t>distribution</w:t></w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r><w:t>normal</w:t></w:r>' +
              '<w:r><w:t>_by_length</w:t></w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '</w:p>'
    $pr1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# Change 2: "preference functions for Categorical layers" heading ->
# "Preference functions for categorical layers" (capitalise P, decapitalise c).
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("preference functions for Categorical layers")
if ($found2) {
    $p2 = $r2.Paragraphs(1)
    $pr2 = $p2.Range
    $xml2 = "<w:p $wNs>" +
              '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
              '<w:r><w:t>P</w:t></w:r>' +
              '<w:r><w:t>reference function</w:t></w:r>' +
              '<w:r><w:t>s</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> for </w:t></w:r>' +
              '<w:r><w:t>c</w:t></w:r>' +
              '<w:r><w:t>ategorical layers</w:t></w:r>' +
            '</w:p>'
    $pr2.InsertXML($xml2)
}


# ---------------------------------------------------------------------------
# Change 3: restructure the "To be further developed..." section.
#  - merge the two "To be further developed and " / "defined in the SPM
#    manual" runs into one run (no more split around the _GoBack bookmark)
#  - move the _GoBack bookmark down to the end of the section (after the
#    last paragraph, replacing the old "Sophie: ..." comment with
#    "No fixed, and tested. It works.")
# ---------------------------------------------------------------------------
$r3a = $d.Content
$found3a = $r3a.Find.Execute("To be further developed and")
$p3a = $r3a.Paragraphs(1)
$start3 = $p3a.Range.Start

$r3b = $d.Content
$found3b = $r3b.Find.Execute("Sophie: Please check with new compiled code")
$p3b = $r3b.Paragraphs(1)
$end3 = $p3b.Range.End

$full3 = $d.Range($start3, $end3)
$xml3 = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>To be further developed and defined in the SPM manual</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Entering layer data</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">At the moment we require a “data” subcommand for every line of SPM to input layer data. Is it worth adopting the form used in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>randomstation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cala</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to have a ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>begin_data</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’ and ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>end’data</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’ style command? Then the subcommand is not required on each new line anymore…?. Note one difference between the programs is that we are inputting a matrix of values by row… and not a table of data with column headings.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Layers</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">At the moment if the data supplied in a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> file has too many data entries, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>is</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> not the correct size (number of rows and columns), SPM fails, the DOS window closes and a Windows error pops up. When layers get loaded, there is a need to check that they are the correct shape (i.e. same size as the base layer) and give an appropriate error message and exit from SPM.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>No fixed, and tested</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> It works.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$full3.InsertXML($xml3)


Write-Output "done"
